$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a dropdown data validation (TRUE/FALSE) on the "enabled" column
$validationRange = $ws.Range("D2:D1048576")
$validationRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$validationRange.Validation.ErrorTitle = "Enabled Error"
$validationRange.Validation.ErrorMessage = "You must choose true or false"

# Convert the boolean TRUE values in column D into literal text "TRUE"
$ws.Range("Z1").Formula = '="TRUE"'
$ws.Range("Z1").Copy()
foreach ($row in 2..7) {
  $ws.Cells.Item($row, 4).PasteSpecial(-4163)
}
$ws.Range("Z1").Clear()

$ws.Range("F8").Select()
